# Deploy updated output folder
# Update the "Metadata" sheet of the CodeSystem spreadsheet:
#  - Title changes from "NG Gender or Sex" to "NG-Imm Gender or Sex CS"
#  - Date changes to the new generation timestamp
#  - "Case Sensitive" value ("true") is cleared
#  - "Hierarchy" value ("is-a") is cleared

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

$ws.Range("B5").Value = "NG-Imm Gender or Sex CS"
$ws.Range("B8").Value = "2025-06-24T09:13:37+01:00"
$ws.Range("B16").Value = ""
$ws.Range("B18").Value = ""
